$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 2.47
$ws.Range("H2").Value = 3
$ws.Range("J2").Value = 3.1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 3.3
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 6.58
$ws.Range("O2").Value = 1.28
$ws.Range("P2").Value = 2.97
$ws.Range("S2").Value = 1.39
$ws.Range("U2").Value = 1.72
$ws.Range("V2").Value = 2.07
$ws.Range("W2").Value = 7.1
$ws.Range("Y2").Value = 7.9
$ws.Range("AA2").Value = 16.5
$ws.Range("AB2").Value = 22
$ws.Range("AC2").Value = 9
$ws.Range("AD2").Value = 5.2
$ws.Range("AE2").Value = 10.5
$ws.Range("AH2").Value = 7.5
$ws.Range("AI2").Value = 11.75
$ws.Range("AK2").Value = 25
$ws.Range("AL2").Value = 18
$ws.Range("AM2").Value = 23
$ws.Range("AN2").Value = 4.4
$ws.Range("AO2").Value = 13.5
$ws.Range("AP2").Value = 21
$ws.Range("AQ2").Value = 60
$ws.Range("AR2").Value = 90
$ws.Range("AS2").Value = 250
$ws.Range("AU2").Value = 6.7
$ws.Range("AV2").Value = 60
$ws.Range("AX2").Value = 15
$ws.Range("AY2").Value = 22
$ws.Range("AZ2").Value = 70
$ws.Range("BA2").Value = 100

# Row 3 updates
$ws.Range("G3").Value = 2.9
$ws.Range("I3").Value = 2.3
$ws.Range("K3").Value = 2.2
$ws.Range("L3").Value = 3
$ws.Range("Q3").Value = 1.95
$ws.Range("R3").Value = 1.85
$ws.Range("Z3").Value = 29
$ws.Range("AB3").Value = 29
$ws.Range("AC3").Value = 11
$ws.Range("AH3").Value = 8.5
$ws.Range("AP3").Value = 23
$ws.Range("AR3").Value = 67
$ws.Range("AX3").Value = 13
